# Add new columns I ("I0") and J ("IF") to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells, styled like the other header cells (style index 1 = bold/border header style)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:J83, one (I,J) pair per data row (rows 2..83)
$data = @(
    @(6,6), @(6,6), @(6,6), @(8,8), @(6,7), @(9,9), @(6,6), @(6,6), @(7,8), @(6,6),
    @(5,5), @(8,8), @(6,6), @(6,7), @(6,7), @(7,7), @(8,8), @(8,8), @(6,6), @(7,7),
    @(9,9), @(6,6), @(7,7), @(8,8), @(6,7), @(9,9), @(6,6), @(8,8), @(8,8), @(10,10),
    @(6,6), @(5,5), @(8,8), @(7,7), @(7,7), @(6,7), @(1,1), @(7,7), @(9,9), @(6,7),
    @(7,7), @(8,8), @(9,9), @(9,9), @(9,9), @(5,6), @(5,6), @(6,6), @(6,7), @(8,8),
    @(7,7), @(8,8), @(1,1), @(7,7), @(7,7), @(7,7), @(9,9), @(8,8), @(7,7), @(5,6),
    @(7,7), @(7,8), @(7,7), @(6,6), @(8,8), @(7,7), @(7,8), @(7,7), @(8,8), @(6,6),
    @(7,8), @(8,8), @(7,8), @(7,8), @(5,5), @(6,6), @(4,5), @(5,6), @(4,4), @(6,6),
    @(7,7), @(5,5)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row = $row + 1
}
